# Update column F ("dSF") values on the active worksheet to reflect the
# repulled/recalculated data from the commit "repull data, push all data,
# mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = -2
    4  = -3
    5  = -4
    6  = 4
    9  = 3
    10 = -1
    12 = 1
    13 = 2
    15 = 4
    17 = -1
    18 = -3
    19 = 1
    20 = 2
    21 = 1
    22 = 2
    23 = 5
    24 = 1
    25 = 3
    27 = 2
    29 = -2
    30 = 1
    31 = -3
    32 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
